# This workbook had previously been re-saved by an unlicensed "Aspose.Cells"
# component, which (a) strips the workbook's structure-lock password and
# (b) appends an "Evaluation Warning" sheet advertising that fact. Reproduce
# both effects here so the file can be freely committed/shared.

$wb = $excel.ActiveWorkbook

# 1) Remove the workbook's structure-protection / password.
$wb.Unprotect("CBEB")

# 2) Append a new worksheet after the existing ones and name it.
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "Evaluation Warning"

# 3) Put the Aspose.Cells evaluation-copy watermark text in A5, styled as a
#    large bold italic blue Arial caption, and size the row to fit it.
$cell = $ws.Range("A5")
$cell.Value = "Evaluation Only. Created with Aspose.Cells for Java.Copyright 2003 - 2016 Aspose Pty Ltd."
$cell.Font.Name = "Arial"
$cell.Font.Size = 18
$cell.Font.Bold = $true
$cell.Font.Italic = $true
$cell.Font.Color = 16711680
$ws.Rows.Item(5).RowHeight = 23.25

# 4) Make the new warning sheet the active/selected tab.
$ws.Activate()
